# Update "想去人数" (number of people interested) figures in column F
# for both the "展览" sheet and the aggregated "全部类型" sheet.
# These two sheets list largely the same events, so the same row/value
# updates are applied to each.

$wb = $excel.ActiveWorkbook

$updates = @{
    5  = 47
    6  = 352
    7  = 10873
    8  = 403
    9  = 92
    11 = 81
    12 = 141
    13 = 143
    20 = 1106
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F" + $row).Value = $updates[$row]
    }
}
